# Auto-generated Excel COM-interop script applying the Marilith_Profits market-data refresh.
# For each changed leve row (per-sheet), update currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) to the refreshed values. Cells that no longer carry a value in the target
# state are cleared (ClearContents) so they serialize as absent, matching the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 49.2
$ws.Range("I38").Value = 49.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 147.6
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 224.4

# Row 61
$ws.Range("H61").Value = 465
$ws.Range("I61").Value = 465
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1395
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1223

# Row 62
$ws.Range("H62").Value = 5432.923
$ws.Range("I62").Value = 4804.4287
$ws.Range("J62").Value = 6166.1665
$ws.Range("K62").Value = 4804.4287
$ws.Range("L62").Value = 6166.1665
$ws.Range("M62").Value = -4180.4287
$ws.Range("N62").Value = -7414.1665

# Row 65
$ws.Range("H65").Value = 5432.923
$ws.Range("I65").Value = 4804.4287
$ws.Range("J65").Value = 6166.1665
$ws.Range("K65").Value = 24022.1435
$ws.Range("L65").Value = 30830.8325
$ws.Range("M65").Value = -20902.1435
$ws.Range("N65").Value = -37070.8325

# Row 141
$ws.Range("H141").Value = 4172.737
$ws.Range("I141").Value = 3963.8235
$ws.Range("J141").Value = 5948.5
$ws.Range("K141").Value = 11891.4705
$ws.Range("L141").Value = 17845.5
$ws.Range("M141").Value = -6711.470499999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5645.184
$ws.Range("I32").Value = 4446.4053
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 4446.4053
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = -4159.4053

# Row 95
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 113
$ws.Range("H113").Value = 21000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 21000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 21000
$ws.Range("N113").Value = -29678

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 95
$ws.Range("H95").Value = 29379
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 29379
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 29379
$ws.Range("N95").Value = -34871

# Row 107
$ws.Range("H107").Value = 1239.5
$ws.Range("I107").Value = 1309.25
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 1309.25
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = 610.75
$ws.Range("N107").Value = -4940

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 5000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -5224

# Row 16
$ws.Range("H16").Value = 2834.8235
$ws.Range("I16").Value = 1506.5714
$ws.Range("J16").Value = 9033.333000000001
$ws.Range("K16").Value = 1506.5714
$ws.Range("L16").Value = 9033.333000000001
$ws.Range("M16").Value = -1219.5714
$ws.Range("N16").Value = -9607.333000000001

# Row 19
$ws.Range("H19").Value = 648.86957
$ws.Range("I19").Value = 376.82352
$ws.Range("J19").Value = 1419.6666
$ws.Range("K19").Value = 376.82352
$ws.Range("L19").Value = 1419.6666
$ws.Range("M19").Value = -206.82352
$ws.Range("N19").Value = -1759.6666

# Row 24
$ws.Range("H24").Value = 648.86957
$ws.Range("I24").Value = 376.82352
$ws.Range("J24").Value = 1419.6666
$ws.Range("K24").Value = 376.82352
$ws.Range("L24").Value = 1419.6666
$ws.Range("M24").Value = -206.82352
$ws.Range("N24").Value = -1759.6666

# Row 42
$ws.Range("H42").Value = 2330.6667
$ws.Range("I42").Value = 2330.6667
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 2330.6667
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -1737.6667
$ws.Range("N42").ClearContents()

# Row 113
$ws.Range("H113").Value = 2834.8235
$ws.Range("I113").Value = 1506.5714
$ws.Range("J113").Value = 9033.333000000001
$ws.Range("K113").Value = 1506.5714
$ws.Range("L113").Value = 9033.333000000001
$ws.Range("M113").Value = 663.4286
$ws.Range("N113").Value = -13373.333

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 836.6667
$ws.Range("I39").Value = 836.6667
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 2510.0001
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -2216.0001
$ws.Range("N39").ClearContents()

# Row 55
$ws.Range("H55").Value = 99
$ws.Range("I55").Value = 99
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 297
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -120
$ws.Range("N55").ClearContents()

# Row 131
$ws.Range("H131").Value = 2154.6
$ws.Range("I131").Value = 1432.3334
$ws.Range("J131").Value = 2636.111
$ws.Range("K131").Value = 4297.0002
$ws.Range("L131").Value = 7908.333
$ws.Range("M131").Value = 742.9997999999996
$ws.Range("N131").Value = -17988.333

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 500
$ws.Range("M5").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 999
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 999
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 999
$ws.Range("N2").Value = -1223

# Row 22
$ws.Range("H22").Value = 278
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 278
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 278
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -868

# Row 27
$ws.Range("H27").Value = 278
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 278
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 278
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -492

# Row 46
$ws.Range("H46").Value = 3464.5881
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 4322
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 4322
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -4698

# Row 55
$ws.Range("H55").Value = 712.9
$ws.Range("I55").Value = 818.4286
$ws.Range("J55").Value = 466.66666
$ws.Range("K55").Value = 818.4286
$ws.Range("L55").Value = 466.66666
$ws.Range("M55").Value = -645.4286

# Row 61
$ws.Range("H61").Value = 1833.1428
$ws.Range("I61").Value = 1833.1428
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1833.1428
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1631.1428

# Row 68
$ws.Range("H68").Value = 6367.25
$ws.Range("I68").Value = 5737.6
$ws.Range("J68").Value = 7416.6665
$ws.Range("K68").Value = 5737.6
$ws.Range("L68").Value = 7416.6665
$ws.Range("M68").Value = -4988.6

# Row 71
$ws.Range("H71").Value = 6367.25
$ws.Range("I71").Value = 5737.6
$ws.Range("J71").Value = 7416.6665
$ws.Range("K71").Value = 28688
$ws.Range("L71").Value = 37083.3325
$ws.Range("M71").Value = -24944

# Row 113
$ws.Range("H113").Value = 1833.1428
$ws.Range("I113").Value = 1833.1428
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1833.1428
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 336.8571999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Range("H74").Value = 21774.6
$ws.Range("I74").Value = 20400.5
$ws.Range("J74").Value = 22690.666
$ws.Range("K74").Value = 20400.5
$ws.Range("L74").Value = 22690.666
$ws.Range("M74").Value = -19464.5
$ws.Range("N74").Value = -24562.666

# Row 77
$ws.Range("H77").Value = 21774.6
$ws.Range("I77").Value = 20400.5
$ws.Range("J77").Value = 22690.666
$ws.Range("K77").Value = 61201.5
$ws.Range("L77").Value = 68071.99800000001
$ws.Range("M77").Value = -56521.5
$ws.Range("N77").Value = -77431.99800000001
